# Updates the cryptos price/volume table (columns D and E, rows 2-51)
# to the refreshed values from the latest scrape.
# Values are written with a leading apostrophe to force text entry
# (many of the "Price" values look numeric, e.g. "1.008"), then the
# cell style is reset back to "Normal" so no visible quote-prefix /
# number-format is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.513.91"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -2.61%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.811.58"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -2.26%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.008"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.70%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'1.007"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.56%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'308.22"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -1.80%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.4555"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -2.03%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.3664"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.07132"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -2.23%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.8794"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -1.29%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07759"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -1.43%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'19.35"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -3.76%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.806.91"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -2.02%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'5.287"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'6.370"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -2.31%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'86.61"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -5.10%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'1.008"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.64%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.000008585"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -3.80%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'1.006"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.54%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'26.594.07"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -2.42%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'14.25"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -3.13%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'5.007"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -1.58%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  -0.57%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'1.983"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +2.11%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'151.53"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +0.09%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  -2.42%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'2.064"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +1.16%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'112.85"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -2.66%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'4.843"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -3.91%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.08679"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -1.74%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'3.063"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -2.45%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'4.523"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +0.13%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.7322"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -4.85%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'2.667"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -1.83%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.118"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -4.23%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'1.005"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +0.63%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'1.083"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -2.52%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.01952"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.50%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05110"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -2.09%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'2.906"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -1.41%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'6.978"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -0.94%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.4998"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -2.49%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  -4.03%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'8.176"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -3.83%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'1.007"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.64%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.4605"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -4.02%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'10.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -3.62%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'101.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -1.67%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'1.592"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -3.21%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  -3.25%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  -1.40%  "
$ws.Range("E51").Style = "Normal"
